$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 214, shifting existing rows 214+ down by one
$ws.Rows.Item(214).Insert()

# Fill in the new row 214 with the new data
$ws.Cells.Item(214, 1).Value = 4
$ws.Cells.Item(214, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(214, 3).Value = "Los Lagos"
$ws.Cells.Item(214, 4).Value = 44830
$ws.Cells.Item(214, 5).Value = 10
$ws.Cells.Item(214, 6).Value = 100112037
$ws.Cells.Item(214, 7).Value = "Cebollín"
$ws.Cells.Item(214, 8).Value = "Sin especificar"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 70
$ws.Cells.Item(214, 11).Value = 9000
$ws.Cells.Item(214, 12).Value = 9000
$ws.Cells.Item(214, 13).Value = 9000
$ws.Cells.Item(214, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(214, 15).Value = "Región Metropolitana"
$ws.Cells.Item(214, 16).Value = 250
$ws.Cells.Item(214, 17).Value = 36
$ws.Cells.Item(214, 18).Value = "Hortaliza"
